$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Metadata" ---
$ws1 = $wb.Worksheets.Item("Metadata")

# Version 5.0.0 -> 6.0.0
$ws1.Range("B3").Value = "6.0.0"

# Date updated
$ws1.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was empty, now "Alvearie Team"
$ws1.Range("B9").Value = "Alvearie Team"

# The second "Contact" / "No display for ContactDetail" row (row 11) is removed;
# row 10 ("Contact" / "No display for ContactDetail") becomes "Jurisdiction" / "United States of America"
$ws1.Range("A10").Value = "Jurisdiction"
$ws1.Range("B10").Value = "United States of America"

# Delete the now-duplicate row 11, shifting rows 12-15 up by one
$ws1.Rows(11).Delete()

# --- Sheet 2: "Include from Claim Care Team " ---
# Content of this sheet is unchanged (only shared-string indices shifted upstream).
